# "6 hours by turn fix"
# Rebuild the Lucas Ferreira weekly schedule sheet: shift the afternoon
# rows so the day now runs through 18:20, inserting an extra time slot,
# and update the subject values that moved between turns/days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The schedule grows from 14 rows (A1:F14) to 17 rows (A1:F17): an extra
# time-slot row is inserted after the current row 13 so everything from
# the old row 14 onward shifts down by one, then two more rows (17:30 and
# 18:20) are appended at the bottom.
$ws.Rows.Item(14).Insert()

# Final values for every row/column of the schedule (time column A, then
# the five weekday columns B:F).
$data = @{
    1  = @("Lucas Ferreira", "segunda", "terça", "quarta", "quinta", "sexta")
    2  = @("7:00", "-", "-", "-", "-", "-")
    3  = @("7:50", "MCT-2A-Sistemas digitais", "-", "-", "-", "-")
    4  = @("8:40", "MCT-2A-Sistemas digitais", "-", "-", "-", "-")
    5  = @("9:30", "Intervalo", "Intervalo", "Intervalo", "Intervalo", "Intervalo")
    6  = @("9:50", "-", "-", "-", "-", "-")
    7  = @("10:40", "-", "-", "ELT-2A-Sistemas digitais", "-", "ELT-2A-Sistemas digitais")
    8  = @("11:30", "-", "-", "-", "-", "-")
    9  = @("12:20", "Almoço", "Almoço", "Almoço", "Almoço", "Almoço")
    10 = @("13:00", "-", "-", "-", "-", "-")
    11 = @("13:50", "-", "-", "-", "-", "-")
    12 = @("14:40", "-", "-", "-", "-", "-")
    13 = @("15:30", "Intervalo", "Intervalo", "Intervalo", "Intervalo", "Intervalo")
    14 = @("15:50", "-", "-", "-", "-", "-")
    15 = @("16:40", "-", "-", "-", "-", "-")
    16 = @("17:30", "-", "-", "-", "-", "-")
    17 = @("18:20", "", "", "", "", "")
}

foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}
